# Refresh cryptos Price (D) and Volume(1h) (E) columns with latest values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.092.60"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "2.589.07"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'524.51"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'139.37"
$ws.Range("E6").Value = "  -3.14%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("D9").Value = "2.600.19"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").Value = "'6.53"
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "3.048.61"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "59.010.48"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").Value = "'20.49"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "2.572.50"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").Value = "'340.79"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("D21").Value = "'10.09"
$ws.Range("E21").Value = "  -1.97%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'66.82"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "'0.406"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "0.0₃0724"
$ws.Range("E30").Value = "  -3.39%  "

$ws.Range("D31").Value = "'5.93"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("D32").Value = "'1.58"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "'18.72"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D34").Value = "'149.20"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("E35").Value = "  -1.37%  "

$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("D37").Value = "'36.80"
$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  +1.42%  "

$ws.Range("D39").Value = "'0.828"
$ws.Range("E39").Value = "  -4.75%  "

$ws.Range("E40").Value = "  -6.42%  "

$ws.Range("E41").Value = "  -0.66%  "

$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").Value = "'271.36"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("D44").Value = "'10.78"
$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("D45").Value = "'0.597"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").Value = "'0.0514"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").Value = "'18.40"
$ws.Range("E48").Value = "  -2.31%  "

$ws.Range("D49").Value = "1.971.49"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "'0.0223"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").Value = "'18.09"
$ws.Range("E51").Value = "  -4.25%  "
